$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Checked/updated base failure rates ---

# Diode base failure rate (left table) and its mirror on the right table
$ws.Range("C3").Value = 0.001965
$ws.Range("J3").Formula = "=C3"

# Mosfet base failure rate (left table)
$ws.Range("C4").Value = 0.01784

# IGBT failure rate (right table)
$ws.Range("J4").Value = 0.015855

# --- Modified converter failure-rate formula (weight by component count) ---
$ws.Range("C8").Formula = "=B2*D2+B3*(D3+D4)+D5+D6"

# New derived cells added around the converter block
$ws.Range("D8").Formula = "=C8*10^6"
$ws.Range("B9").Formula = "=EXP(-C8*0.2*E4)"
$ws.Range("C10").Formula = "=EXP(-C8*E6*0.2)"
$ws.Range("J10").Formula = "=J9*10^6"

# --- Passive balancing table ---
$ws.Range("C20").Value = 0.003966
$ws.Range("C21").Formula = "=C4"
$ws.Range("E24").Formula = "=0.2*E22"

# --- Column width for new column J content ---
$ws.Columns.Item(10).ColumnWidth = 11.25

# --- Selection moved to J11 ---
$ws.Range("J11").Select()
